$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("K_AGE", "Altersklasse", "Age group"),
    @("K_AIRPOLL", "Luftschadstoff", "Air pollutant"),
    @("K_AREA", "Gebiet", "Area"),
    @("K_CALCMETH", "Berechnungsmethode", "Calculation method"),
    @("K_CRIM", "Straftat", "Criminal offence"),
    @("K_CRIMOFF", "Straftaten", "Criminal offences"),
    @("K_KREIS", "Kreis", "County"),
    @("K_LAENDER", "Bundesland", "Federal state"),
    @("K_PM", "Feinstaub", "Fine particulate matter"),
    @("K_QUALI", "Qualifizierung", "XXX"),
    @("K_SEA", "Meer", "Sea"),
    @("K_SERIES", "Zeitreihe", "Time series"),
    @("K_SEX", "Geschlecht", "Sex"),
    @("K_SUBINDEX", "Teilindizes", "Sub index"),
    @("K_TYPEAREA", "Art der Fläche", "Type of area"),
    @("K_URBAN", "Verstädterungsgrad", "Degree of urbanisation")
)

# Copy the formatting of the existing data row (row 2) so the newly
# inserted rows (3-17) pick up the same style (s="4") before we
# overwrite row 2's own values further below.
$formatSource = $ws.Range("A2:C2")
$formatSource.Copy()
$ws.Range("A3:C17").PasteSpecial(-4122)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
